$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- simple literal replacements (dates / amounts / acordadas) ---
Replace-All "01/02/2023" "21/02/2021"
Replace-All "`$6.853.038,49" "`$3.306.575,22"
Replace-All "27/06/2023" "22/05/2022"
Replace-All "Acordada 19/2023" "Acordada 12/2022"
Replace-All "`$19.338,00" "`$9.001,00"
Replace-All "354.38" "367.36"
Replace-All "16/04/2024" "25/07/2024"
Replace-All "`$13.229.256,23" "`$13.562.727,77"
Replace-All "Acordada 1497/2024" "Acordada 1772/2024"
Replace-All "`$52.510,00" "`$57.016,00"
Replace-All "130.51" "237.88"

"done"
